# Update the "Förändrad" (Changed) date column (C) for all data rows:
# the stored serial date 45204 (2023-10-05) becomes 45205 (2023-10-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (xlUp = -4162), falling back to
# the worksheet's UsedRange if that lookup fails for some reason.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
